$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.647.33'
$ws.Range("E2").Value = '  -4.32%  '

$ws.Range("D3").Value = '3.340.52'
$ws.Range("E3").Value = '  -1.34%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''574.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.31%  '

$ws.Range("D6").Value = '''180.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.22%  '

$ws.Range("D7").Value = '''0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.33%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -3.25%  '

$ws.Range("D10").Value = '''6.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.76%  '

$ws.Range("D11").Value = '''0.403'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = '3.918.40'
$ws.Range("E12").Value = '  -1.45%  '

$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").Value = '''27.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.77%  '

$ws.Range("D15").Value = '66.746.52'
$ws.Range("E15").Value = '  -4.10%  '

$ws.Range("D16").Value = '''0.0000167'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").Value = '3.350.06'
$ws.Range("E17").Value = '  -1.09%  '

$ws.Range("D18").Value = '''436.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.50%  '

$ws.Range("E19").Value = '  -2.36%  '

$ws.Range("D20").Value = '''13.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").Value = '''7.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.82%  '

$ws.Range("D22").Value = '''73.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.42%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = '''0.517'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("E25").Value = '  -3.96%  '

$ws.Range("E26").Value = '  -0.48%  '

$ws.Range("D27").Value = '''9.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.01%  '

$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("E29").Value = '  -3.40%  '

$ws.Range("D30").Value = '''22.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.69%  '

$ws.Range("E31").Value = '  -6.30%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").Value = '''6.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.22%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '''1.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.26%  '

$ws.Range("D35").Value = '''162.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.38%  '

$ws.Range("E36").Value = '  -5.89%  '

$ws.Range("D37").Value = '''27.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("D38").Value = '''1.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.09%  '

$ws.Range("D39").Value = '2.825.72'
$ws.Range("E39").Value = '  +2.57%  '

$ws.Range("D40").Value = '''0.795'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.46%  '

$ws.Range("E41").Value = '  -3.96%  '

$ws.Range("E42").Value = '  -6.05%  '

$ws.Range("E43").Value = '  -2.33%  '

$ws.Range("E44").Value = '  -3.19%  '

$ws.Range("D45").Value = '''24.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.42%  '

$ws.Range("E46").Value = '  -6.45%  '

$ws.Range("D47").Value = '''321.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.45%  '

$ws.Range("E48").Value = '  -4.02%  '

$ws.Range("E49").Value = '  +1.00%  '

$ws.Range("E50").Value = '  -4.08%  '

$ws.Range("D51").Value = '''6.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.57%  '
